$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: insert a new sub-bullet "Remove it from the Datalab
# preparation" (ilvl=1) right after the "Go into correlation
# coefficient..." bullet and before "Assessment Criteria..." bullet.
# -----------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Go into correlation coefficient more deeply, get more clear understanding ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Collapse(0)
$rng1.InsertAfter("`rRemove it from the Datalab preparation")

$rng1b = $d.Content
$found1b = $rng1b.Find.Execute("Remove it from the Datalab preparation", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1b.ListFormat.ListLevelNumber = 2
$rng1b.LanguageID = "en-NL"

# -----------------------------------------------------------------
# Edit 2: append five new bullets (ilvl=0) at the very end of the
# document, after "Define appropriate more clearly: conditions to
# use boxplot/histogram etc." -- the last one is an empty bullet.
# -----------------------------------------------------------------
$rng2 = $d.Content
$rng2.Collapse(0)
$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Calculation </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t>excercises</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t xml:space="preserve"> need more explanations; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t>expecially</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t xml:space="preserve"> specifying the symbols and terms in the equations and what they mean</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t>Sigma calculation</w:t></w:r><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t xml:space="preserve"> put in again</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t>Group students together who work on the same SDG indicators; sessions where students work together on their SDG indicator separate of their mentor group</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t xml:space="preserve">Consider replacing the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t>datacamp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t xml:space="preserve"> courses with other materials; especially less courses on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t>datacamp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t xml:space="preserve"> later on and connect it with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t>excercises</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-NL"/></w:rPr><w:t xml:space="preserve"> which mirror the dashboard deliverables</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($xmlFrag)
